$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet view: zoom scale and selection
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("N10").Select() | Out-Null

# Row 2
$ws.Range("B2").Value = 100
$ws.Range("C2").Formula = "=B2*3"
$ws.Range("F2").Value = 0.4
$ws.Range("L2").Value = 0.147

# Row 3
$ws.Range("B3").Value = 200
$ws.Range("C3").Formula = "=B3*3"
$ws.Range("F3").Value = 0.4
$ws.Range("L3").Value = 0.147

# Row 4
$ws.Range("B4").Value = 500
$ws.Range("C4").Formula = "=B4*3"
$ws.Range("F4").Value = 0.4
$ws.Range("L4").Value = 0.147

# Row 5
$ws.Range("B5").Value = 2500
$ws.Range("C5").Formula = "=B5*3"
$ws.Range("F5").Value = 0.4
$ws.Range("L5").Value = 0.147

# Row 6 (B6 unchanged)
$ws.Range("C6").Formula = "=B6*3"
$ws.Range("F6").Value = 0.4
$ws.Range("L6").Value = 0.147
